$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update masses (column D) for rows 6-10; dependent ixx/iyy/izz formulas
# in E:G recalculate automatically. ---
$ws.Range("D6").Value = 4
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 0.1
$ws.Range("D9").Value = 0.05
$ws.Range("D10").Value = 0.05

# --- Give the inertia results for the two finger links more decimal
# precision so the small values are legible. ---
$ws.Range("E9:G10").NumberFormat = "0.0000000"

# --- Resize the result columns to fit their new (longer) content. ---
$ws.Range("E1:G10").EntireColumn.AutoFit()

# --- Switch the printable page to portrait orientation. ---
$ws.PageSetup.Orientation = 1

# --- Restore the active cell selection. ---
$ws.Range("I15").Select() | Out-Null
